$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = rotation speed -> set to 1 for rows 2-7
$ws.Range("A2:A7").Value = 1

# Column B = mask speed -> set to 2 for rows 2-7, except row 4 which is 1
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 2

# Update the current selection to B7 as in the saved file
$ws.Range("B7").Select()
